$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 430
$ws.Range("I18").Value = 430
$ws.Range("K18").Value = 430
$ws.Range("M18").Value = -146

$ws.Range("H28").Value = 319
$ws.Range("I28").Value = 142.9
$ws.Range("J28").Value = 1199.5
$ws.Range("K28").Value = 142.9
$ws.Range("L28").Value = 1199.5
$ws.Range("M28").Value = 342.1
$ws.Range("N28").Value = -2169.5

$ws.Range("H69").Value = 1546.6
$ws.Range("I69").Value = 1874.5
$ws.Range("J69").Value = 1496.1538
$ws.Range("K69").Value = 5623.5
$ws.Range("L69").Value = 4488.4614
$ws.Range("M69").Value = -4749.5
$ws.Range("N69").Value = -6236.4614

$ws.Range("H72").Value = 1546.6
$ws.Range("I72").Value = 1874.5
$ws.Range("J72").Value = 1496.1538
$ws.Range("K72").Value = 16870.5
$ws.Range("L72").Value = 13465.3842
$ws.Range("M72").Value = -12502.5
$ws.Range("N72").Value = -22201.3842

$ws.Range("H88").Value = 5004
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5004
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5004
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -5816

$ws.Range("H91").Value = 5004
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5004
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5004
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -7812

$ws.Range("H129").Value = 162379.16
$ws.Range("I129").Value = 497.42856
$ws.Range("J129").Value = 182982.3
$ws.Range("K129").Value = 1492.28568
$ws.Range("L129").Value = 548946.8999999999
$ws.Range("M129").Value = 3507.71432
$ws.Range("N129").Value = -558946.8999999999

$ws.Range("H132").Value = 2393.8809
$ws.Range("I132").Value = 2572.0264
$ws.Range("J132").Value = 701.5
$ws.Range("K132").Value = 7716.0792
$ws.Range("L132").Value = 2104.5
$ws.Range("M132").Value = -5186.0792
$ws.Range("N132").Value = -7164.5

$ws.Range("H135").Value = 13893436
$ws.Range("I135").Value = 280
$ws.Range("K135").Value = 2520
$ws.Range("M135").Value = 15

$ws.Range("H141").Value = 1254.7826
$ws.Range("I141").Value = 1130
$ws.Range("K141").Value = 3390
$ws.Range("M141").Value = 1790

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4898.913
$ws.Range("I32").Value = 5299.65
$ws.Range("K32").Value = 5299.65
$ws.Range("M32").Value = -5012.65

$ws.Range("H74").Value = 25642962
$ws.Range("I74").Value = 32259892
$ws.Range("K74").Value = 32259892
$ws.Range("M74").Value = -32259018

$ws.Range("H77").Value = 25642962
$ws.Range("I77").Value = 32259892
$ws.Range("K77").Value = 161299460
$ws.Range("M77").Value = -161295092

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3320.1538
$ws.Range("I20").Value = 4079.2
$ws.Range("J20").Value = 790
$ws.Range("K20").Value = 4079.2
$ws.Range("L20").Value = 790
$ws.Range("M20").Value = -3832.2
$ws.Range("N20").Value = -1284

$ws.Range("H134").Value = 3200.7441
$ws.Range("I134").Value = 3577.7646
$ws.Range("J134").Value = 1776.4445
$ws.Range("K134").Value = 10733.2938
$ws.Range("L134").Value = 5329.333500000001
$ws.Range("M134").Value = -8198.293799999999
$ws.Range("N134").Value = -10399.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 247.5
$ws.Range("I22").Value = 188.42857
$ws.Range("J22").Value = 330.2
$ws.Range("K22").Value = 188.42857
$ws.Range("L22").Value = 330.2
$ws.Range("M22").Value = 161.57143
$ws.Range("N22").Value = -1030.2

$ws.Range("H31").Value = 16723.545
$ws.Range("I31").Value = 30990.4
$ws.Range("J31").Value = 4834.5
$ws.Range("K31").Value = 30990.4
$ws.Range("L31").Value = 4834.5
$ws.Range("M31").Value = -30695.4
$ws.Range("N31").Value = -5424.5

$ws.Range("H34").Value = 16723.545
$ws.Range("I34").Value = 30990.4
$ws.Range("J34").Value = 4834.5
$ws.Range("K34").Value = 30990.4
$ws.Range("L34").Value = 4834.5
$ws.Range("M34").Value = -30788.4
$ws.Range("N34").Value = -5238.5

$ws.Range("H53").Value = 35634.5
$ws.Range("J53").Value = 35634.5
$ws.Range("L53").Value = 35634.5
$ws.Range("N53").Value = -36848.5

$ws.Range("H58").Value = 21234.4
$ws.Range("I58").Value = 1263.8125
$ws.Range("K58").Value = 1263.8125
$ws.Range("M58").Value = -1060.8125

$ws.Range("H69").Value = 9068.5
$ws.Range("J69").Value = 9191.333000000001
$ws.Range("L69").Value = 9191.333000000001
$ws.Range("N69").Value = -10689.333

$ws.Range("H72").Value = 9068.5
$ws.Range("J72").Value = 9191.333000000001
$ws.Range("L72").Value = 27573.999
$ws.Range("N72").Value = -35061.999

$ws.Range("H132").Value = 16394.143
$ws.Range("I132").Value = 19312.143
$ws.Range("K132").Value = 57936.429
$ws.Range("M132").Value = -55406.429

$ws.Range("H134").Value = 940.8261
$ws.Range("I134").Value = 823.1053000000001
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 2469.3159
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = 65.68409999999994
$ws.Range("N134").Value = -9570

$ws.Range("H136").Value = 21234.4
$ws.Range("I136").Value = 1263.8125
$ws.Range("K136").Value = 3791.4375
$ws.Range("M136").Value = -1241.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 72.5
$ws.Range("J12").Value = 98.454544
$ws.Range("L12").Value = 295.363632
$ws.Range("N12").Value = -641.3636320000001

$ws.Range("H117").Value = 3280.5
$ws.Range("I117").Value = 1929
$ws.Range("J117").Value = 4632
$ws.Range("K117").Value = 5787
$ws.Range("L117").Value = 13896
$ws.Range("M117").Value = -2345
$ws.Range("N117").Value = -20780

$ws.Range("H122").Value = 441.6
$ws.Range("I122").Value = 250.23077
$ws.Range("K122").Value = 2252.07693
$ws.Range("M122").Value = 197.9230699999998

$ws.Range("H131").Value = 122774.95
$ws.Range("I131").Value = 827.8
$ws.Range("J131").Value = 130693.59
$ws.Range("K131").Value = 2483.4
$ws.Range("L131").Value = 392080.77
$ws.Range("M131").Value = 2556.6
$ws.Range("N131").Value = -402160.77

$ws.Range("H138").Value = 188920.5
$ws.Range("I138").Value = 1451.1111
$ws.Range("K138").Value = 4353.3333
$ws.Range("M138").Value = 786.6666999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1973.3334
$ws.Range("I68").Value = 1628.5714
$ws.Range("J68").Value = 2275
$ws.Range("K68").Value = 1628.5714
$ws.Range("L68").Value = 2275
$ws.Range("M68").Value = -879.5714
$ws.Range("N68").Value = -3773

$ws.Range("H71").Value = 1973.3334
$ws.Range("I71").Value = 1628.5714
$ws.Range("J71").Value = 2275
$ws.Range("K71").Value = 8142.857
$ws.Range("L71").Value = 11375
$ws.Range("M71").Value = -4398.857
$ws.Range("N71").Value = -18863

$ws.Range("H132").Value = 1559.2354
$ws.Range("I132").Value = 1385.1538
$ws.Range("K132").Value = 4155.4614
$ws.Range("M132").Value = -1625.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 166668740
$ws.Range("I81").Value = 2480
$ws.Range("K81").Value = 4960
$ws.Range("M81").Value = -3899

$ws.Range("H84").Value = 166668740
$ws.Range("I84").Value = 2480
$ws.Range("K84").Value = 24800
$ws.Range("M84").Value = -19496

$ws.Range("H132").Value = 1675.05
$ws.Range("I132").Value = 1060.3
$ws.Range("J132").Value = 2289.8
$ws.Range("K132").Value = 3180.9
$ws.Range("L132").Value = 6869.400000000001
$ws.Range("M132").Value = -650.8999999999996
$ws.Range("N132").Value = -11929.4

$ws.Range("H136").Value = 38463236
$ws.Range("I136").Value = 45456096
$ws.Range("J136").Value = 2499.75
$ws.Range("K136").Value = 136368288
$ws.Range("L136").Value = 7499.25
$ws.Range("M136").Value = -136365738
